$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

# Drop the second sheet - only one sheet remains in the target workbook.
[void]$wb.Worksheets.Item("mySheet_2").Delete()

# Rename the remaining sheet and make sure it is the active one.
$ws = $wb.Worksheets.Item("Blatt1")
$ws.Name = "blatt1"
[void]$ws.Activate()

# Clear whatever sample data was there before.
[void]$ws.Cells.Clear()

# Header / misc string
$ws.Range("B1").Value = "This is another test string"

# Summary block
$ws.Range("D2").Value = "Zeitraum:"
$ws.Range("E2").NumberFormat = "mm-dd-yy"
$ws.Range("E2").Value = Get-Date -Year 2017 -Month 1 -Day 1 -Hour 0 -Minute 0 -Second 0
$ws.Range("D3").Value = "Summe [km]:"
$ws.Range("D4").Value = "Kilometergeld:"

# Table header row
$ws.Range("A6").Value = "start_generated"
$ws.Range("B6").Value = "Datum"
$ws.Range("C6").Value = "km-Stand"
$ws.Range("D6").Value = "Start"
$ws.Range("E6").Value = "Ziel"
$ws.Range("F6").Value = "km"
$ws.Range("G6").Value = "Kommentar"

# End marker
$ws.Range("A8").Value = "end_generated"

[void]$ws.Range("A9").Select()
